{"js": "// Add a new row to the (only) table with the grants to-do item, then\n// append a few new paragraphs after the table.\n\nconst body = context.document.body;\n\n// --- 1. Add the new table row -------------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.addRows(Word.InsertLocation.end, 1, [\n  [\n    \"Add up our total grants since the start of the pandemic and compare to the total grants from the NIPAs since the start of the pandemic to determine if we\\u2019re on target\",\n    \"8/17/2023\",\n    \"To do\"\n  ]\n]);\nawait context.sync();\n\n// --- 2. Append the new paragraphs after the table ------------------------------\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"Calculate grants\", Word.InsertLocation.end);\nbody.insertParagraph(\n  \"Figure out why transfers change as a result of grants\",\n  Word.InsertLocation.end\n);\nbody.insertParagraph(\n  \"When we move forward a quarter, make sure we \",\n  Word.InsertLocation.end\n);\nawait context.sync();\n", "ps1": "# Add a new row to the (only) table with the grants to-do item, then\n# append a few new paragraphs after the table.\n\n$d = $word.ActiveDocument\n\n# --- 1. Add the new table row -------------------------------------------------\n$table = $d.Tables.Item(1)\n$newRow = $table.Rows.Add()\n$idx = $newRow.Index\n\n$table.Cell($idx, 1).Range.Text = \"Add up our total grants since the start of the pandemic and compare to the total grants from the NIPAs since the start of the pandemic to determine if we\" + [char]0x2019 + \"re on target\"\n$table.Cell($idx, 2).Range.Text = \"8/17/2023\"\n$table.Cell($idx, 3).Range.Text = \"To do\"\n\n# --- 2. Append the new paragraphs after the table ------------------------------\n# NB: re-derive the end-of-document range from $d.Content.End each time\n# instead of caching $d.Paragraphs.Last/.Count - re-using a Paragraphs\n# collection snapshot taken before the table row insert above can resolve\n# to stale indices once the table has grown.\n\n# empty paragraph\n$d.Range($d.Content.End - 1, $d.Content.End - 1).InsertParagraphAfter()\n\n# \"Calculate grants\"\n$d.Range($d.Content.End - 1, $d.Content.End - 1).InsertParagraphAfter()\n$d.Range($d.Content.End - 1, $d.Content.End - 1).InsertAfter(\"Calculate grants\")\n\n# \"Figure out why transfers change as a result of grants\"\n$d.Range($d.Content.End - 1, $d.Content.End - 1).InsertParagraphAfter()\n$d.Range($d.Content.End - 1, $d.Content.End - 1).InsertAfter(\"Figure out why transfers change as a result of grants\")\n\n# \"When we move forward a quarter, make sure we \"\n$d.Range($d.Content.End - 1, $d.Content.End - 1).InsertParagraphAfter()\n$d.Range($d.Content.End - 1, $d.Content.End - 1).InsertAfter(\"When we move forward a quarter, make sure we \")\n"}
